$d = $word.ActiveDocument

# 1. "Your child has just taken a STAR " -> "Your child has taken a STAR "
$d.Content.Find.Execute("Your child has just taken a STAR ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Your child has taken a STAR ", 2)

# 2. "...your child's scores on the assessment. As with any assessment," ->
#    "...your child's scores. As with any assessment,"
# (the search/replace text below deliberately starts after the apostrophe in
#  "child's" so the replacement string itself contains no apostrophe -- Word's
#  smart-quote autocorrect otherwise turns the existing straight apostrophe
#  into a curly one whenever the replacement text is "typed")
$d.Content.Find.Execute("s scores on the assessment. As with any assessment,", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "s scores. As with any assessment,", 2)

# 3. "in each of the four domains; Numbers" -> "in each of the four domains: Numbers"
$d.Content.Find.Execute("in each of the four domains; Numbers", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "in each of the four domains: Numbers", 2)

# 4. "You are the first teacher to your child. To help " -> "You are the first teacher for your child. To help "
$d.Content.Find.Execute("You are the first teacher to your child. To help ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "You are the first teacher for your child. To help ", 2)

# 5. " mathematical skills try doing some of these activities at home. " ->
#    " mathematical skills try doing some of these activities at home: "
$d.Content.Find.Execute(" mathematical skills try doing some of these activities at home. ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " mathematical skills try doing some of these activities at home: ", 2)

# 6. Move the "_GoBack" bookmark: it used to sit at the end of the opening
#    paragraph (right after "...of how your child is doing in school."); it
#    now belongs right after "...activities at home:" (before the trailing
#    space) further down the letter. Adding a bookmark named "_GoBack" removes
#    any pre-existing "_GoBack" bookmark elsewhere in the document, since Word
#    only ever keeps a single one.
$marker = $d.Content
$marker.Find.Execute("ome of these activities at home:", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0)
$bmRange = $d.Range($marker.End, $marker.End)
$d.Bookmarks.Add("_GoBack", $bmRange)
